# Rename worksheets (tab names) to new "summ" identifiers, keeping order,
# sheetId and r:id unchanged - only the <sheet name="..."> values change.

$wb = $excel.ActiveWorkbook

$renames = @(
    @{ Old = "summ17839566"; New = "summ06142548" },
    @{ Old = "summ18076052"; New = "summ06386333" },
    @{ Old = "summ18326490"; New = "summ06647780" },
    @{ Old = "summ18578178"; New = "summ06884356" },
    @{ Old = "summ18830578"; New = "summ07135129" },
    @{ Old = "summ19102323"; New = "summ07408808" },
    @{ Old = "summ19365234"; New = "summ07653564" },
    @{ Old = "summ19621403"; New = "summ07902504" },
    @{ Old = "summ19876324"; New = "summ08172851" }
)

foreach ($r in $renames) {
    $sheet = $wb.Worksheets.Item($r.Old)
    $sheet.Name = $r.New
}
